$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 4744
$ws.Range("F5").Value = 206
$ws.Range("F6").Value = 1906
$ws.Range("F8").Value = 795
$ws.Range("F11").Value = 422
$ws.Range("F12").Value = 1163
$ws.Range("F14").Value = 839
$ws.Range("F15").Value = 33
$ws.Range("F16").Value = 1937
$ws.Range("F17").Value = 598
$ws.Range("F18").Value = 37
$ws.Range("F19").Value = 536
$ws.Range("F20").Value = 630
$ws.Range("F21").Value = 226
$ws.Range("F22").Value = 78
$ws.Range("F23").Value = 78
$ws.Range("F26").Value = 626
$ws.Range("F27").Value = 2538
$ws.Range("F28").Value = 21
$ws.Range("F31").Value = 1636
$ws.Range("F33").Value = 504
$ws.Range("F36").Value = 4349

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 4171
$ws.Range("F11").Value = 45
$ws.Range("F16").Value = 295
$ws.Range("F17").Value = 295
$ws.Range("F21").Value = 149
$ws.Range("F23").Value = 250
$ws.Range("F29").Value = 108
$ws.Range("F38").Value = 42

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1369
$ws.Range("F5").Value = 1751
$ws.Range("F7").Value = 386

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1369
$ws.Range("F3").Value = 1751
$ws.Range("F5").Value = 386
$ws.Range("F9").Value = 4744
$ws.Range("F11").Value = 206
$ws.Range("F12").Value = 1906
$ws.Range("F14").Value = 795
$ws.Range("F19").Value = 422
$ws.Range("F20").Value = 1163
$ws.Range("F23").Value = 839
$ws.Range("F24").Value = 36
$ws.Range("F25").Value = 1937
$ws.Range("F26").Value = 598
$ws.Range("F27").Value = 40
$ws.Range("F28").Value = 536
$ws.Range("F29").Value = 226
$ws.Range("F30").Value = 78
$ws.Range("F31").Value = 295
$ws.Range("F32").Value = 295
$ws.Range("F36").Value = 626
$ws.Range("F37").Value = 149
$ws.Range("F39").Value = 2538
$ws.Range("F40").Value = 250
$ws.Range("F45").Value = 1636
$ws.Range("F46").Value = 504
$ws.Range("F49").Value = 4349
$ws.Range("F50").Value = 42
